$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.786.61'
$ws.Range("E2").Value = '  -0.13%  '

$ws.Range("D3").Value = '2.825.85'
$ws.Range("E3").Value = '  +1.24%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '350.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.10%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '112.84'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.59%  '

$ws.Range("E7").Value = '  +1.38%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("E9").Value = '  +3.71%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.18'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.72%  '

$ws.Range("E11").Value = '  -0.90%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0848'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.22%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.07'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.58%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.79'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.56%  '

$ws.Range("D15").Value = '3.271.04'
$ws.Range("E15").Value = '  +1.39%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.983'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.02%  '

$ws.Range("D17").Value = '2.818.28'
$ws.Range("E17").Value = '  +2.17%  '

$ws.Range("D18").Value = '51.787.38'
$ws.Range("E18").Value = '  +0.20%  '

$ws.Range("E19").Value = '  +9.48%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.63'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.82%  '

$ws.Range("E21").Value = '  +1.50%  '

$ws.Range("E22").Value = '  +0.65%  '

$ws.Range("E23").Value = '  +0.47%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '268.65'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.54%  '

$ws.Range("E25").Value = '  +0.57%  '

$ws.Range("E26").Value = '  +0.18%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.06%  '

$ws.Range("E28").Value = '  +0.79%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '39.17'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +7.41%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.52'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.47%  '

$ws.Range("E31").Value = '  +6.52%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.33'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.27%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '52.72'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.28%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0895'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.70%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.63'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.52%  '

$ws.Range("E36").Value = '  -0.97%  '

$ws.Range("E37").Value = '  +0.34%  '

$ws.Range("E38").Value = '  +1.53%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.22'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.89%  '

$ws.Range("E40").Value = '  +1.86%  '

$ws.Range("E41").Value = '  +1.07%  '

$ws.Range("E42").Value = '  -1.56%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '122.93'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.83%  '

$ws.Range("E44").Value = '  +1.49%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '22.09'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.12%  '

$ws.Range("E46").Value = '  +8.54%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.51'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.75%  '

$ws.Range("D48").Value = '2.168.36'
$ws.Range("E48").Value = '  +2.32%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.248'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +22.03%  '

$ws.Range("E50").Value = '  +3.04%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.49'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.59%  '
